# Update the "想去人数" (interest count) figures in column F across the
# relevant worksheets, matching the data refresh performed by the
# gh-pages output generation (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 49
$ws1.Range("F4").Value = 312
$ws1.Range("F5").Value = 1262
$ws1.Range("F7").Value = 300
$ws1.Range("F8").Value = 1098
$ws1.Range("F9").Value = 428
$ws1.Range("F10").Value = 6901
$ws1.Range("F14").Value = 7803
$ws1.Range("F17").Value = 5184
$ws1.Range("F19").Value = 2285
$ws1.Range("F20").Value = 969
$ws1.Range("F21").Value = 4539
$ws1.Range("F22").Value = 253
$ws1.Range("F23").Value = 367
$ws1.Range("F26").Value = 285
$ws1.Range("F29").Value = 2017
$ws1.Range("F31").Value = 226
$ws1.Range("F33").Value = 531
$ws1.Range("F35").Value = 1377
$ws1.Range("F36").Value = 21
$ws1.Range("F37").Value = 2101

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 89

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 49
$ws4.Range("F7").Value = 312
$ws4.Range("F8").Value = 1262
$ws4.Range("F11").Value = 300
$ws4.Range("F12").Value = 1098
$ws4.Range("F13").Value = 428
$ws4.Range("F14").Value = 6901
$ws4.Range("F18").Value = 7803
$ws4.Range("F21").Value = 5184
$ws4.Range("F23").Value = 2285
$ws4.Range("F24").Value = 969
$ws4.Range("F25").Value = 4539
$ws4.Range("F26").Value = 253
$ws4.Range("F27").Value = 367
$ws4.Range("F32").Value = 285
$ws4.Range("F35").Value = 2017
$ws4.Range("F37").Value = 226
$ws4.Range("F39").Value = 531
$ws4.Range("F42").Value = 1377
$ws4.Range("F43").Value = 21
$ws4.Range("F44").Value = 2101
$ws4.Range("F48").Value = 89
